# Updates the "Price" (column D) and "Volume(1h)" (column E) figures
# for the crypto-ranking snapshot, matching the refreshed scrape.
# Values are written with a leading apostrophe so Excel keeps them as
# plain text (matching the original inlineStr cells) instead of
# auto-converting numeric-looking / percentage-looking strings into
# actual numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'329.92"
$ws.Range("E2").Value = "'7.28%"
# Row 3
$ws.Range("D3").Value = "'40.24"
$ws.Range("E3").Value = "'8.90%"
# Row 4
$ws.Range("D4").Value = "'5.282"
$ws.Range("E4").Value = "'3.26%"
# Row 5
$ws.Range("D5").Value = "'0.08103"
$ws.Range("E5").Value = "'3.42%"
# Row 6
$ws.Range("D6").Value = "'4.512"
$ws.Range("E6").Value = "'2.81%"
# Row 7
$ws.Range("D7").Value = "'8.638"
$ws.Range("E7").Value = "'4.49%"
# Row 8
$ws.Range("D8").Value = "'1.918"
$ws.Range("E8").Value = "'1.48%"
# Row 9
$ws.Range("E9").Value = "'-0.47%"
# Row 10
$ws.Range("D10").Value = "'0.9366"
$ws.Range("E10").Value = "'1.29%"
# Row 11
$ws.Range("D11").Value = "'0.1343"
$ws.Range("E11").Value = "'24.21%"
# Row 12
$ws.Range("D12").Value = "'0.1969"
$ws.Range("E12").Value = "'3.79%"
# Row 13
$ws.Range("D13").Value = "'0.09236"
$ws.Range("E13").Value = "'3.85%"
# Row 14
$ws.Range("D14").Value = "'0.03571"
$ws.Range("E14").Value = "'7.51%"
# Row 15
$ws.Range("D15").Value = "'0.09586"
$ws.Range("E15").Value = "'0.13%"
# Row 16
$ws.Range("D16").Value = "'0.001327"
$ws.Range("E16").Value = "'-3.82%"
# Row 17
$ws.Range("D17").Value = "'0.006126"
$ws.Range("E17").Value = "'3.34%"
# Row 18
$ws.Range("D18").Value = "'3.377"
$ws.Range("E18").Value = "'-3.45%"
# Row 19
$ws.Range("D19").Value = "'0.3523"
$ws.Range("E19").Value = "'3.12%"
# Row 20
$ws.Range("D20").Value = "'7.165"
$ws.Range("E20").Value = "'13.82%"
# Row 21
$ws.Range("D21").Value = "'0.1322"
$ws.Range("E21").Value = "'3.56%"
# Row 22
$ws.Range("E22").Value = "'2.20%"
# Row 23
$ws.Range("D23").Value = "'0.04430"
$ws.Range("E23").Value = "'2.18%"
# Row 24
$ws.Range("D24").Value = "'0.001222"
$ws.Range("E24").Value = "'2.54%"
# Row 25
$ws.Range("D25").Value = "'0.004350"
$ws.Range("E25").Value = "'2.34%"
# Row 26
$ws.Range("E26").Value = "'-8.57%"
# Row 27
$ws.Range("D27").Value = "'0.0003990"
$ws.Range("E27").Value = "'-0.04%"
# Row 39
$ws.Range("D39").Value = "'0.02501"
$ws.Range("E39").Value = "'16.42%"
# Row 40
$ws.Range("D40").Value = "'0.05185"
$ws.Range("E40").Value = "'3.37%"
# Row 41
$ws.Range("D41").Value = "'0.007664"
$ws.Range("E41").Value = "'2.51%"
# Row 43
$ws.Range("D43").Value = "'0.009132"
$ws.Range("E43").Value = "'5.50%"
# Row 44
$ws.Range("D44").Value = "'0.002170"
$ws.Range("E44").Value = "'2.64%"
# Row 45
$ws.Range("D45").Value = "'0.01008"
$ws.Range("E45").Value = "'14.52%"
# Row 46
$ws.Range("D46").Value = "'0.00006651"
$ws.Range("E46").Value = "'1.43%"
# Row 47
$ws.Range("E47").Value = "'-0.13%"
# Row 48
$ws.Range("D48").Value = "'0.002400"
$ws.Range("E48").Value = "'139.69%"
# Row 49
$ws.Range("E49").Value = "'17.22%"
# Row 50
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.13%"
# Row 51
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.13%"
